$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Grafico de burndown" backlog item entirely (row 24) - the team
# decided to add presentation slides instead of a burndown chart task.
$ws.Rows.Item(24).Delete()

# Mark the remaining previously-open items as finished ("Finalizado"),
# since the whole backlog is now complete.
$ws.Range("I15").Value = "Finalizado"
$ws.Range("I16").Value = "Finalizado"
$ws.Range("I17").Value = "Finalizado"
$ws.Range("I21").Value = "Finalizado"
$ws.Range("I22").Value = "Finalizado"
$ws.Range("I23").Value = "Finalizado"
$ws.Range("I24").Value = "Finalizado"
$ws.Range("I25").Value = "Finalizado"

# Shrink conditional formatting ranges to follow the new (smaller) data range
$dFcs = $ws.Range("D3:D25").FormatConditions()
for ($i = 1; $i -le $dFcs.Count(); $i++) {
  $dFcs.Item($i).ModifyAppliesToRange($ws.Range("D3:D25"))
}
$iFcs = $ws.Range("I3:I25").FormatConditions()
for ($i = 1; $i -le $iFcs.Count(); $i++) {
  $iFcs.Item($i).ModifyAppliesToRange($ws.Range("I3:I25"))
}
